$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ B = 4949.642832416669;  D = 204.6528945666667 }
    3  = @{ B = 4654.605469183336;  D = 187.6830644333333 }
    4  = @{ B = 4945.102591533336;  D = 209.5629226166667 }
    5  = @{ B = 4605.854491083336;  D = 193.0902597666667 }
    6  = @{ B = 4934.395209500003;  D = 199.3819656833333 }
    7  = @{ B = 4645.800626650002;  D = 194.8143918833333 }
    8  = @{ B = 4874.360841583336;  D = 208.1396406833333 }
    9  = @{ B = 4869.838935383335;  D = 191.7460731333333 }
    10 = @{ B = 4605.612003033336;  D = 190.79789355 }
    11 = @{ B = 4706.058441483336;  D = 191.4939004166667 }
    12 = @{ B = 4787.523504650003;  D = 202.5828063666667 }
    13 = @{ B = 4743.635858583336;  D = 197.0560429 }
}

foreach ($row in $newValues.Keys) {
    $ws.Range("B$row").Value = $newValues[$row].B
    $ws.Range("D$row").Value = $newValues[$row].D
}

$wb.Save()
